$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell from "Communes" to "commune"
$ws.Range("B1").Value = "commune"

# Move the current selection to C24 (single cell), matching the saved view state
$ws.Range("C24").Select()
